# Update "want to go" (想去人数) counts in column F on the
# "展览" (Exhibition), "演出" (Performance) and "全部类型" (All types) sheets.
# Values taken from the authoritative diff of the generated output.
# NOTE: this COM runtime does not reliably bind *named* PowerShell
# parameters, so the helper function below uses positional parameters.

$wb = $excel.ActiveWorkbook

function Set-FValues($SheetName, $RowValues) {
    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($row in $RowValues.Keys) {
        $ws.Cells.Item($row, 6).Value = $RowValues[$row]
    }
}

# 展览 sheet (sheet1)
Set-FValues "展览" @{
    5  = 15745
    6  = 417
    8  = 708
    9  = 15459
    11 = 9060
    12 = 388
    14 = 1015
    16 = 202
    21 = 558
    25 = 1116
    29 = 94
    31 = 42
    34 = 44
    39 = 5569
}

# 演出 sheet (sheet2)
Set-FValues "演出" @{
    2 = 69
}

# 全部类型 sheet (sheet4)
Set-FValues "全部类型" @{
    5  = 15745
    6  = 417
    8  = 708
    9  = 15459
    11 = 9060
    12 = 388
    14 = 1015
    16 = 202
    21 = 558
    25 = 1116
    29 = 94
    31 = 42
    32 = 69
    36 = 44
    41 = 5569
}

$wb.Save()
